# Update "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, per the commit "Update gh-pages to output generated
# at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 10789
$ws.Range("F4").Value = 258
$ws.Range("F5").Value = 1196
$ws.Range("F6").Value = 1065
$ws.Range("F7").Value = 836
$ws.Range("F8").Value = 277
$ws.Range("F9").Value = 341
$ws.Range("F10").Value = 1142
$ws.Range("F12").Value = 875
$ws.Range("F14").Value = 1925
$ws.Range("F16").Value = 949
$ws.Range("F17").Value = 818
$ws.Range("F19").Value = 793
$ws.Range("F20").Value = 904
$ws.Range("F25").Value = 628
$ws.Range("F26").Value = 116
$ws.Range("F28").Value = 1006
$ws.Range("F31").Value = 162
$ws.Range("F33").Value = 230
$ws.Range("F34").Value = 556
$ws.Range("F35").Value = 1770
$ws.Range("F36").Value = 374
$ws.Range("F38").Value = 1421
$ws.Range("F39").Value = 408
$ws.Range("F40").Value = 122
$ws.Range("F41").Value = 48
$ws.Range("F42").Value = 83
$ws.Range("F45").Value = 73
$ws.Range("F46").Value = 79

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 1
$ws.Range("F10").Value = 29
$ws.Range("F14").Value = 129

# Sheet: 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2165
$ws.Range("F3").Value = 625
$ws.Range("F4").Value = 554

# Sheet: 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2165
$ws.Range("F3").Value = 625
$ws.Range("F5").Value = 10789
$ws.Range("F6").Value = 258
$ws.Range("F8").Value = 554
$ws.Range("F9").Value = 1065
$ws.Range("F10").Value = 836
$ws.Range("F12").Value = 277
$ws.Range("F13").Value = 1142
$ws.Range("F18").Value = 1925
$ws.Range("F21").Value = 793
$ws.Range("F22").Value = 904
$ws.Range("F27").Value = 29
$ws.Range("F28").Value = 628
$ws.Range("F29").Value = 116
$ws.Range("F31").Value = 1006
$ws.Range("F33").Value = 162
$ws.Range("F34").Value = 230
$ws.Range("F36").Value = 129
$ws.Range("F37").Value = 374
$ws.Range("F38").Value = 1421
$ws.Range("F39").Value = 408
$ws.Range("F40").Value = 122
$ws.Range("F42").Value = 48
$ws.Range("F43").Value = 83
$ws.Range("F45").Value = 73
